# Commit: "Fixed product name typo & adding tiny images"
#
# 1) Product "butterfly-earings" (row 5 of Main / row 4 of Categories) had a
#    typo in its slug, name and description -> "butterfly-earrings".
# 2) Product "brown-fedora" (row 4 of Main) was using a mismatched/placeholder
#    image ("indiana-jones-hat.jpg") -> switched to its own tiny image
#    ("brown-fedora-01.jpg").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Main")

# Fix the "Butterfly Earings" -> "Butterfly Earrings" typo (slug, name, description)
$ws1.Range("A5").Value = "butterfly-earrings"
$ws1.Range("E5").Value = "Butterfly Earrings"
$ws1.Range("M5").Value = "Sample Butterfly Earrings Lorem ipsum dolor sit amet, consectetur adipisicing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum"

# Swap in the correct (tiny) product image for the Brown Fedora row
$ws1.Range("L4").Value = "brown-fedora-01.jpg"

# Same slug fix on the "Categories" sheet (row for butterfly-earrings)
$ws2 = $wb.Worksheets.Item("Categories")
$ws2.Range("A4").Value = "butterfly-earrings"

# Restore view/selection state: the selected cell on "Categories" moves to A9,
# while "Main" stays the active tab with L5 selected.
$ws2.Range("A9").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("L5").Select() | Out-Null
